# Generate Report for Handoff
# Updates "Latest Handoff Datetime" (column D) for rows whose handoff file
# was just (re-)generated, on both the "zh-cn" and "de-de" status sheets.
#
# Row 4  -> bd494f5a-a353-4424-8e66-0a6f21fd7afa
# Row 6  -> 5e371f2c-00a3-4c4e-a8e2-4b34a60323e9
# Row 7  -> 807d69c9-0db5-4f6f-baa6-87306a5461e1
# Row 8  -> 98297fad-c161-441e-856d-82c8fbbaac22
# Row 9  -> 9c37f105-1d40-4858-a041-69277b88b4ac
# Row 10 -> e1600472-db79-4b30-be19-18c3c021dff9
#
# Rows 2, 3, 5 and 11 are left untouched (already in sync / in translation /
# not localized).

$wb = $excel.ActiveWorkbook

$rows = @(4, 6, 7, 8, 9, 10)

$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 4).Value = "2016-02-22 18:19:07"
}

$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 4).Value = "2016-02-22 18:19:21"
}
